$d = $word.ActiveDocument

# 1. Replace the French text with the template placeholder {o1}
$d.Content.Find.Execute("Véritable copie de l'original", $true, $false, $false, $false, $false, $true, 1, $false, "{o1}", 2)

# 2. Move the "_GoBack" bookmark so that it now sits right after the
#    newly-inserted "{o1}" run (it used to live in the previous, empty
#    paragraph). Locate the paragraph that now contains "{o1}".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{o1}*") {
        $target = $p
    }
}

$endPos = $target.Range.End - 1

# NOTE: placing a *collapsed* bookmark exactly one character before a
# paragraph mark can land in the wrong spot, so we temporarily insert a
# marker character after the insertion point, add/move the bookmark there
# (now safely not adjacent to the paragraph mark), and then remove the
# marker again. The collapsed bookmark stays put when the marker is
# deleted.
$insertRng = $d.Range($endPos, $endPos)
$insertRng.InsertAfter("Z")

$bmRng = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

$delRng = $d.Range($endPos, $endPos + 1)
$delRng.Delete()
